$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dSF = @{
    2 = 0
    3 = -3
    6 = 5
    7 = -3
    8 = 1
    9 = -6
    10 = 2
    11 = -3
    12 = -2
    13 = 2
    14 = -1
    15 = 3
    16 = 2
    17 = -3
    19 = -3
    20 = 1
    21 = 2
    22 = -3
    23 = 4
    24 = 2
    25 = -3
    26 = -1
    27 = -1
    28 = -1
    29 = -3
    30 = -6
    31 = 4
    32 = 3
    33 = 1
    34 = -3
    35 = -3
}

foreach ($row in $dSF.Keys) {
    $ws.Cells.Item($row, 6).Value = $dSF[$row]
}
